# COREESG_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer (A10) from
#    2021-04-06 to 2021-04-08
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-6
#    and the Percent Change total (E7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; drop protection so the cells can be written,
# then restore it (same password used in the source workbook) afterwards.
$ws.Unprotect("D382")

$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-08 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.2512921186425509
$ws.Range("E2").Value = 0.0115423218467714

$ws.Range("D3").Value = 0.4930744432138315
$ws.Range("E3").Value = -0.0008178844056706547

$ws.Range("D4").Value = 0.09983538050603893
$ws.Range("E4").Value = 0.01811805961426072

$ws.Range("D5").Value = 0.09906885764842009
$ws.Range("E5").Value = -0.002573634543894698

$ws.Range("D6").Value = 0.05672919998915867
$ws.Range("E6").Value = 0.007395424081349855

$ws.Range("E7").Value = 0.00447060944605826

$ws.Protect("D382")
